# Auto-generated cell updates for rows 207-214 (Poland Ekstraklasa odds refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 207
$ws.Cells.Item(207, 2).Value = 6775564
$ws.Cells.Item(207, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(207, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(207, 5).Value = 45359.58333333334
$ws.Cells.Item(207, 6).Value = "Piast Gliwice"
$ws.Cells.Item(207, 7).Value = "Radomiak Radom"
$ws.Cells.Item(207, 8).Value = 2
$ws.Cells.Item(207, 9).Value = 3
$ws.Cells.Item(207, 10).Value = "A"
$ws.Cells.Item(207, 11).Value = 2
$ws.Cells.Item(207, 12).Value = 3.25
$ws.Cells.Item(207, 13).Value = 4
$ws.Cells.Item(207, 14).Value = 2
$ws.Cells.Item(207, 15).Value = 3
$ws.Cells.Item(207, 16).Value = 4.333
$ws.Cells.Item(207, 17).Value = -0.5
$ws.Cells.Item(207, 18).Value = 2.025
$ws.Cells.Item(207, 19).Value = 1.825
$ws.Cells.Item(207, 20).Value = 2
$ws.Cells.Item(207, 21).Value = 2.05
$ws.Cells.Item(207, 22).Value = 1.8
$ws.Cells.Item(207, 23).Value = -1
$ws.Cells.Item(207, 24).Value = -1
$ws.Cells.Item(207, 25).Value = 3.333
$ws.Cells.Item(207, 26).Value = -1
$ws.Cells.Item(207, 27).Value = 0.825
$ws.Cells.Item(207, 28).Value = 1.05
$ws.Cells.Item(207, 29).Value = -1

# Row 208
$ws.Cells.Item(208, 2).Value = 6775562
$ws.Cells.Item(208, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(208, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(208, 5).Value = 45359.6875
$ws.Cells.Item(208, 6).Value = "Jagiellonia Bialystok"
$ws.Cells.Item(208, 7).Value = "Slask Wroclaw"
$ws.Cells.Item(208, 8).Value = 3
$ws.Cells.Item(208, 9).Value = 1
$ws.Cells.Item(208, 10).Value = "H"
$ws.Cells.Item(208, 11).Value = 1.952
$ws.Cells.Item(208, 12).Value = 3.5
$ws.Cells.Item(208, 13).Value = 3.8
$ws.Cells.Item(208, 14).Value = 1.95
$ws.Cells.Item(208, 15).Value = 3.4
$ws.Cells.Item(208, 16).Value = 3.8
$ws.Cells.Item(208, 17).Value = -0.5
$ws.Cells.Item(208, 18).Value = 2
$ws.Cells.Item(208, 19).Value = 1.85
$ws.Cells.Item(208, 20).Value = 2.5
$ws.Cells.Item(208, 21).Value = 2
$ws.Cells.Item(208, 22).Value = 1.85
$ws.Cells.Item(208, 23).Value = 0.95
$ws.Cells.Item(208, 24).Value = -1
$ws.Cells.Item(208, 25).Value = -1
$ws.Cells.Item(208, 26).Value = 1
$ws.Cells.Item(208, 27).Value = -1
$ws.Cells.Item(208, 28).Value = 1
$ws.Cells.Item(208, 29).Value = -1

# Row 209
$ws.Cells.Item(209, 2).Value = 6774879
$ws.Cells.Item(209, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(209, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(209, 5).Value = 45360.5625
$ws.Cells.Item(209, 6).Value = "MKS Puszcza Niepolomice"
$ws.Cells.Item(209, 7).Value = "Rakow Czestochowa"
$ws.Cells.Item(209, 11).Value = 4.75
$ws.Cells.Item(209, 12).Value = 3.8
$ws.Cells.Item(209, 13).Value = 1.7
$ws.Cells.Item(209, 14).Value = 5.5
$ws.Cells.Item(209, 15).Value = 4
$ws.Cells.Item(209, 16).Value = 1.571
$ws.Cells.Item(209, 17).Value = 1
$ws.Cells.Item(209, 18).Value = 1.825
$ws.Cells.Item(209, 19).Value = 2.025
$ws.Cells.Item(209, 20).Value = 2.5
$ws.Cells.Item(209, 21).Value = 1.875
$ws.Cells.Item(209, 22).Value = 1.975
$ws.Cells.Item(209, 23).Value = 0
$ws.Cells.Item(209, 24).Value = 0
$ws.Cells.Item(209, 25).Value = 0
$ws.Cells.Item(209, 26).Value = 0
$ws.Cells.Item(209, 27).Value = 0

# Row 210
$ws.Cells.Item(210, 2).Value = 6775561
$ws.Cells.Item(210, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(210, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(210, 5).Value = 45360.66666666666
$ws.Cells.Item(210, 6).Value = "Gornik Zabrze"
$ws.Cells.Item(210, 7).Value = "Lech Poznan"
$ws.Cells.Item(210, 11).Value = 3
$ws.Cells.Item(210, 12).Value = 3.4
$ws.Cells.Item(210, 13).Value = 2.3
$ws.Cells.Item(210, 14).Value = 3
$ws.Cells.Item(210, 15).Value = 3.3
$ws.Cells.Item(210, 16).Value = 2.375
$ws.Cells.Item(210, 17).Value = 0.25
$ws.Cells.Item(210, 18).Value = 1.775
$ws.Cells.Item(210, 19).Value = 2.1
$ws.Cells.Item(210, 20).Value = 2.25
$ws.Cells.Item(210, 21).Value = 1.8
$ws.Cells.Item(210, 22).Value = 2.05
$ws.Cells.Item(210, 23).Value = 0
$ws.Cells.Item(210, 24).Value = 0
$ws.Cells.Item(210, 25).Value = 0
$ws.Cells.Item(210, 26).Value = 0
$ws.Cells.Item(210, 27).Value = 0

# Row 211
$ws.Cells.Item(211, 2).Value = 6774464
$ws.Cells.Item(211, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(211, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(211, 5).Value = 45361.35416666666
$ws.Cells.Item(211, 6).Value = "Stal Mielec"
$ws.Cells.Item(211, 7).Value = "Ruch Chorzow"
$ws.Cells.Item(211, 11).Value = 2
$ws.Cells.Item(211, 12).Value = 3.4
$ws.Cells.Item(211, 13).Value = 3.75
$ws.Cells.Item(211, 14).Value = 2.7
$ws.Cells.Item(211, 15).Value = 3.1
$ws.Cells.Item(211, 16).Value = 2.7
$ws.Cells.Item(211, 17).Value = 0
$ws.Cells.Item(211, 18).Value = 1.925
$ws.Cells.Item(211, 19).Value = 1.925
$ws.Cells.Item(211, 20).Value = 2.25
$ws.Cells.Item(211, 21).Value = 1.975
$ws.Cells.Item(211, 22).Value = 1.875
$ws.Cells.Item(211, 23).Value = 0
$ws.Cells.Item(211, 24).Value = 0
$ws.Cells.Item(211, 25).Value = 0
$ws.Cells.Item(211, 26).Value = 0
$ws.Cells.Item(211, 27).Value = 0

# Row 212
$ws.Cells.Item(212, 2).Value = 6775565
$ws.Cells.Item(212, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(212, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(212, 5).Value = 45361.45833333334
$ws.Cells.Item(212, 6).Value = "Pogon Szczecin"
$ws.Cells.Item(212, 7).Value = "Zaglebie Lubin"
$ws.Cells.Item(212, 11).Value = 1.666
$ws.Cells.Item(212, 12).Value = 4
$ws.Cells.Item(212, 13).Value = 4.75
$ws.Cells.Item(212, 14).Value = 1.727
$ws.Cells.Item(212, 15).Value = 3.8
$ws.Cells.Item(212, 16).Value = 4.5
$ws.Cells.Item(212, 17).Value = -0.75
$ws.Cells.Item(212, 18).Value = 1.975
$ws.Cells.Item(212, 19).Value = 1.875
$ws.Cells.Item(212, 20).Value = 2.75
$ws.Cells.Item(212, 21).Value = 1.875
$ws.Cells.Item(212, 22).Value = 1.975
$ws.Cells.Item(212, 23).Value = 0
$ws.Cells.Item(212, 24).Value = 0
$ws.Cells.Item(212, 25).Value = 0
$ws.Cells.Item(212, 26).Value = 0
$ws.Cells.Item(212, 27).Value = 0

# Row 213
$ws.Cells.Item(213, 2).Value = 6775566
$ws.Cells.Item(213, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(213, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(213, 5).Value = 45361.5625
$ws.Cells.Item(213, 6).Value = "Widzew Lodz"
$ws.Cells.Item(213, 7).Value = "Legia Warsaw"
$ws.Cells.Item(213, 11).Value = 3.5
$ws.Cells.Item(213, 12).Value = 3.5
$ws.Cells.Item(213, 13).Value = 2.05
$ws.Cells.Item(213, 14).Value = 3.8
$ws.Cells.Item(213, 15).Value = 3.5
$ws.Cells.Item(213, 16).Value = 1.95
$ws.Cells.Item(213, 17).Value = 0.5
$ws.Cells.Item(213, 18).Value = 1.85
$ws.Cells.Item(213, 19).Value = 2
$ws.Cells.Item(213, 20).Value = 2.5
$ws.Cells.Item(213, 21).Value = 1.95
$ws.Cells.Item(213, 22).Value = 1.9
$ws.Cells.Item(213, 23).Value = 0
$ws.Cells.Item(213, 24).Value = 0
$ws.Cells.Item(213, 25).Value = 0
$ws.Cells.Item(213, 26).Value = 0
$ws.Cells.Item(213, 27).Value = 0

# Row 214
$ws.Cells.Item(214, 2).Value = 6774465
$ws.Cells.Item(214, 3).Value = "Poland Ekstraklasa"
$ws.Cells.Item(214, 4).Value = "Poland Ekstraklasa"
$ws.Cells.Item(214, 5).Value = 45362.625
$ws.Cells.Item(214, 6).Value = "Warta Poznan"
$ws.Cells.Item(214, 7).Value = "LKS Lodz"
$ws.Cells.Item(214, 11).Value = 1.95
$ws.Cells.Item(214, 12).Value = 3.3
$ws.Cells.Item(214, 13).Value = 4
$ws.Cells.Item(214, 14).Value = 1.95
$ws.Cells.Item(214, 15).Value = 3.3
$ws.Cells.Item(214, 16).Value = 4
$ws.Cells.Item(214, 17).Value = -0.5
$ws.Cells.Item(214, 18).Value = 1.95
$ws.Cells.Item(214, 19).Value = 1.9
$ws.Cells.Item(214, 20).Value = 2.25
$ws.Cells.Item(214, 21).Value = 2
$ws.Cells.Item(214, 22).Value = 1.85
$ws.Cells.Item(214, 23).Value = 0
$ws.Cells.Item(214, 24).Value = 0
$ws.Cells.Item(214, 25).Value = 0
$ws.Cells.Item(214, 26).Value = 0
$ws.Cells.Item(214, 27).Value = 0

# Remove the now-obsolete last row (old row 215), shifting nothing else
$ws.Rows(215).Delete()
